$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison": update MyForecast (column D) values ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws1.Range("D8").Value = 82
$ws1.Range("D9").Value = 64
$ws1.Range("D10").Value = 58
$ws1.Range("D11").Value = 62
$ws1.Range("D13").Value = 76
$ws1.Range("D15").Value = 68

# --- Sheet "Summary": update derived totals to match the refreshed forecast ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B10").Value = "799"
$ws2.Range("B14").Value = "58"
